# Generate Report for Archive
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it is used (Overview!E2:F2, and the
#   per-language sheets' Status column, C2).
# - Shrink the now-narrower Status columns (Overview E:F, zh-cn C, de-de C)
#   to reflect the shorter text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status text -------------------------------------------------
if ($overview.Range("E2").Text -eq $oldStatus) { $overview.Range("E2").Value = $newStatus }
if ($overview.Range("F2").Text -eq $oldStatus) { $overview.Range("F2").Value = $newStatus }
if ($zhcn.Range("C2").Text -eq $oldStatus) { $zhcn.Range("C2").Value = $newStatus }
if ($dede.Range("C2").Text -eq $oldStatus) { $dede.Range("C2").Value = $newStatus }

# --- Resize the columns that held the status text ---------------------------
$overview.Columns("E:E").ColumnWidth = 12.5
$overview.Columns("F:F").ColumnWidth = 12.5
$zhcn.Columns("C:C").ColumnWidth = 12.5
$dede.Columns("C:C").ColumnWidth = 12.5
